$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: CmsWork
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CmsWork")

$ws.Range("C2").Value = "_:N1b67ed7ad3d74a6985b6de1610cb99d7"
$ws.Range("E2").Value = "CmsCollection0CmsWork1 alternative title 1"
$ws.Range("G2").Value = "CmsCollection0CmsWork1Id1"
$ws.Range("O2").Value = "http://www.wikidata.org/entity/Q937690"
$ws.Range("T2").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:101"

$ws.Range("C3").Value = "_:N9682d3b7fd254a77a2ee9b7e55c9c9e5"
$ws.Range("E3").Value = "CmsCollection0CmsWork3 alternative title 0"
$ws.Range("O3").Value = "http://www.wikidata.org/entity/Q937690"
$ws.Range("T3").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:104"

$ws.Range("C4").Value = "_:N1ecc13cc932542b38b02e061e5e98f5a"
$ws.Range("E4").Value = "CmsCollection1CmsWork5 alternative title 1"
$ws.Range("T4").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:106"

$ws.Range("C5").Value = "_:Nfaaa2930c33d46a18eb51d89e16822e3"
$ws.Range("D5").Value = "http://example.com/person4"
$ws.Range("J5").Value = "CmsCollection1CmsWork7 provenance 0"
$ws.Range("O5").Value = "http://www.wikidata.org/entity/Q937690"

$ws.Range("B6").Value = "_:N3b8beb311c79491aa36df60bf2393e70"
$ws.Range("C6").Value = "http://example.com/organization1"
$ws.Range("D6").Value = "FreestandingWork9 alternative title 0"
$ws.Range("F6").Value = "FreestandingWork9Id1"
$ws.Range("I6").Value = "FreestandingWork9 provenance 0"

$ws.Range("B7").Value = "_:Nbdf8f714e5f14a4bb9ed15b441088b5b"
$ws.Range("F7").Value = "FreestandingWork11Id0"
$ws.Range("I7").Value = "FreestandingWork11 provenance 0"
$ws.Range("N7").Value = "http://www.wikidata.org/entity/Q937690"

# ---------------------------------------------------------------------------
# Sheet: CmsWorkClosing
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CmsWorkClosing")

$ws.Range("A2").Value = "_:N78f9fb8b601b4dd8b244387a47fd655e"
$ws.Range("C2").Value = "_:N3f2c1769594a4d76bd141a9d07512e81"

$ws.Range("A3").Value = "_:N9de5b4a7a8774c339267d1314709a5e7"
$ws.Range("C3").Value = "_:Naebc8b6f28cc4f8aa989f31190202631"

$ws.Range("A4").Value = "_:N8ac5352da322401f9aefdb9752cb99e3"
$ws.Range("C4").Value = "_:N102adedbd7334ec1a0a5191426911adf"

$ws.Range("A5").Value = "_:Nac31e4d88060489ebb7309897f69c567"
$ws.Range("C5").Value = "_:N151e370315f94c979009a2f553ee75de"

$ws.Range("A6").Value = "_:Na7912610025a4598b65b43b70441490c"
$ws.Range("C6").Value = "_:Nb1e1d5606b664221a07f7ca5717c1638"

$ws.Range("A7").Value = "_:Nacc0614712b0414794b225b9feb472cc"
$ws.Range("C7").Value = "_:N755686b71e2949a2b6b1dc97cef8208a"

# ---------------------------------------------------------------------------
# Sheet: CmsWorkOpening
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CmsWorkOpening")

$ws.Range("C2").Value = "_:N3f2c1769594a4d76bd141a9d07512e81"
$ws.Range("C3").Value = "_:Naebc8b6f28cc4f8aa989f31190202631"
$ws.Range("C4").Value = "_:N102adedbd7334ec1a0a5191426911adf"
$ws.Range("C5").Value = "_:N151e370315f94c979009a2f553ee75de"
$ws.Range("C6").Value = "_:Nb1e1d5606b664221a07f7ca5717c1638"
$ws.Range("C7").Value = "_:N755686b71e2949a2b6b1dc97cef8208a"

# ---------------------------------------------------------------------------
# Sheet: CmsRightsStatement
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CmsRightsStatement")

$ws.Range("E2").Value = "You may need to obtain other permissions for your intended use. For example, other rights such as publicity, privacy or moral rights may limit how you may use the material."

# ---------------------------------------------------------------------------
# Sheet: CmsPerson
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CmsPerson")

$ws.Range("F2").Value = "http://en.wikipedia.org/wiki/Alan_Turing"
$ws.Range("E5").Value = "http://en.wikipedia.org/wiki/Alan_Turing"
